{"js": "// Replace the legal-authority placeholder and the \"System Owner\" references\n// in the AP-1 \"Privacy\" paragraph with the real legislation cite and the\n// actual org name, per commit \"added legislation in Privacy controls\".\n\nconst body = context.document.body;\n\n// 1) \"Under , System Owner ...\" -> \"Under Title II of the Workforce\n//    Innovation and Opportunity Act (WIOA), System Owner ...\"\nconst underResults = body.search(\"Under ,\", { matchCase: true, matchWholeWord: false });\nunderResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < underResults.items.length; i++) {\n  underResults.items[i].insertText(\n    \"Under Title II of the Workforce Innovation and Opportunity Act (WIOA),\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 2) Every \"System Owner\" -> \"Example Org\" (two occurrences in this paragraph).\nconst ownerResults = body.search(\"System Owner\", { matchCase: true, matchWholeWord: false });\nownerResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < ownerResults.items.length; i++) {\n  ownerResults.items[i].insertText(\"Example Org\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the legal-authority placeholder and the \"System Owner\" references\n# in the AP-1 \"Privacy\" paragraph with the real legislation cite and the\n# actual org name, per commit \"added legislation in Privacy controls\".\n\n$d = $word.ActiveDocument\n\n# 1) \"Under , System Owner ...\" -> \"Under Title II of the Workforce\n#    Innovation and Opportunity Act (WIOA), System Owner ...\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Under ,\"\n$find1.Replacement.Text = \"Under Title II of the Workforce Innovation and Opportunity Act (WIOA),\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) Every \"System Owner\" -> \"Example Org\" (two occurrences in this paragraph).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"System Owner\"\n$find2.Replacement.Text = \"Example Org\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
